$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Usernames
$ws.Range("B2").Value = "userA"
$ws.Range("B3").Value = "userB"
$ws.Range("B4").Value = "userC"

# Address column - now alphanumeric text values instead of plain zip numbers
$ws.Range("D2").Value = "Adr123"
$ws.Range("D3").Value = "Adr124"
$ws.Range("D4").Value = "Adr125"

# Phone column - format as text ("@") so leading zeros are preserved, then
# write the phone numbers as text strings
$ws.Range("H1:H4").NumberFormat = "@"
$ws.Range("H2").Value = "0712345678"
$ws.Range("H3").Value = "0712345679"
$ws.Range("H4").Value = "0712345680"

# Move the active selection to B4
$ws.Range("B4").Select()
